$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7; existing rows 7-17 shift down to 8-18.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with this week's record (same structure/values as the
# other rows, matching column layout A..R).
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = "2021-10-13"
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112013
$ws.Range("G7").Value = "Alcachofa"
$ws.Range("H7").Value = "Madrigal"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11500
$ws.Range("N7").Value = "$/caja 40 unidades"
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 288
$ws.Range("Q7").Value = 40
$ws.Range("R7").Value = "Hortaliza"
